# Insert a new "MinU18" column into the grade sheet, between "MinFemales"
# (existing column M) and "TeamSize" (existing column N), shifting
# TeamSize/Limit/StartLimit one column to the right, and populate the
# value for the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N (14) currently holds "TeamSize"; inserting here shifts it (and
# the columns after it) one place to the right, making room for the new
# "MinU18" column.
$ws.Columns.Item(14).Insert()

# New column header and value.
$ws.Cells.Item(1, 14).Value = "MinU18"
$ws.Cells.Item(2, 14).Value = 2

# Give the new column a width close to its neighbours instead of the
# sheet default.
$ws.Columns.Item(14).ColumnWidth = 8.83

# Leave the selection on the newly added column, like it was left after
# performing the insert interactively.
$ws.Range("N1").Select() | Out-Null
